# Update "想去人数" (F column) and "最低票价" (G column) values across
# the four worksheets to reflect the latest scrape output.

$wb = $excel.ActiveWorkbook

# --- Sheet: 展览 ---
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 354
$ws.Range("F4").Value = 409
$ws.Range("F5").Value = 1116
$ws.Range("F8").Value = 845
$ws.Range("F9").Value = 1591
$ws.Range("F10").Value = 6031
$ws.Range("F12").Value = 1723
$ws.Range("F13").Value = 436
$ws.Range("F14").Value = 5823
$ws.Range("F15").Value = 111
$ws.Range("F16").Value = 48
$ws.Range("F19").Value = 1628
$ws.Range("F22").Value = 137
$ws.Range("F23").Value = 1329
$ws.Range("G23").Value = 69
$ws.Range("F24").Value = 713
$ws.Range("F25").Value = 225
$ws.Range("F28").Value = 25
$ws.Range("F30").Value = 3851

# --- Sheet: 演出 ---
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F4").Value = 301
$ws.Range("F5").Value = 152
$ws.Range("F8").Value = 370

# --- Sheet: 本地生活 ---
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 9487
$ws.Range("F3").Value = 2213
$ws.Range("F4").Value = 596
$ws.Range("F5").Value = 157

# --- Sheet: 全部类型 ---
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 9487
$ws.Range("F3").Value = 2213
$ws.Range("F4").Value = 596
$ws.Range("F5").Value = 354
$ws.Range("F6").Value = 409
$ws.Range("F7").Value = 1116
$ws.Range("F11").Value = 301
$ws.Range("F12").Value = 846
$ws.Range("F13").Value = 157
$ws.Range("F14").Value = 1591
$ws.Range("F15").Value = 6031
$ws.Range("F17").Value = 1723
$ws.Range("F20").Value = 436
$ws.Range("F23").Value = 5823
$ws.Range("F24").Value = 111
$ws.Range("F25").Value = 48
$ws.Range("F28").Value = 1628
$ws.Range("F31").Value = 137
$ws.Range("F32").Value = 1329
$ws.Range("G32").Value = 69
$ws.Range("F33").Value = 713
$ws.Range("F34").Value = 225
$ws.Range("F40").Value = 25
$ws.Range("F45").Value = 3851
